# Weekly update: insert two new data rows (week of 2021-11-08, Excel serial 44508)
# for "Coliflor" at "Mercado Mayorista Lo Valledor de Santiago", pushing the
# existing rows 390..485 down to 392..487 and growing the sheet from
# A1:R485 to A1:R487.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 390 (shifts everything
# from row 390 downward by two rows, exactly like Excel's own Insert Rows).
$ws.Range("A390:A391").EntireRow.Insert()

# --- New row 390: Coliflor, "Primera" quality --------------------------------
$ws.Cells.Item(390, 1).Value  = 6
$ws.Cells.Item(390, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(390, 3).Value  = "Metropolitana"
$ws.Cells.Item(390, 4).Value  = 44508
$ws.Cells.Item(390, 5).Value  = 13
$ws.Cells.Item(390, 6).Value  = 100112008
$ws.Cells.Item(390, 7).Value  = "Coliflor"
$ws.Cells.Item(390, 8).Value  = "Sin especificar"
$ws.Cells.Item(390, 9).Value  = "Primera"
$ws.Cells.Item(390, 10).Value = 8300
$ws.Cells.Item(390, 11).Value = 600
$ws.Cells.Item(390, 12).Value = 650
$ws.Cells.Item(390, 13).Value = 621
$ws.Cells.Item(390, 14).Value = "$/unidad"
$ws.Cells.Item(390, 15).Value = "Región Metropolitana"
$ws.Cells.Item(390, 16).Value = 621
$ws.Cells.Item(390, 17).Value = 1
$ws.Cells.Item(390, 18).Value = "Hortaliza"

# --- New row 391: Coliflor, "Segunda" quality --------------------------------
$ws.Cells.Item(391, 1).Value  = 6
$ws.Cells.Item(391, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(391, 3).Value  = "Metropolitana"
$ws.Cells.Item(391, 4).Value  = 44508
$ws.Cells.Item(391, 5).Value  = 13
$ws.Cells.Item(391, 6).Value  = 100112008
$ws.Cells.Item(391, 7).Value  = "Coliflor"
$ws.Cells.Item(391, 8).Value  = "Sin especificar"
$ws.Cells.Item(391, 9).Value  = "Segunda"
$ws.Cells.Item(391, 10).Value = 2200
$ws.Cells.Item(391, 11).Value = 450
$ws.Cells.Item(391, 12).Value = 450
$ws.Cells.Item(391, 13).Value = 450
$ws.Cells.Item(391, 14).Value = "$/unidad"
$ws.Cells.Item(391, 15).Value = "Región Metropolitana"
$ws.Cells.Item(391, 16).Value = 450
$ws.Cells.Item(391, 17).Value = 1
$ws.Cells.Item(391, 18).Value = "Hortaliza"
